$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Title cell (A1): drop the "(токойлор)" qualifier from the indicator name
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "15.b.1.1 Айлана-чөйрөнү коргоо үчүн бөлүнгөн мамлекеттик акча каражаттарынын көлөмү"

# ---------------------------------------------------------------------------
# 2) Unit caption (C2, "(million soms)"): pick up the same italic style
#    already used by the Kyrgyz/Russian captions in A2/B2
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Add the 2023 data column (Q). First clone the number formatting from
#    column O (already using the trimmed "0.0" format) across both the
#    existing column P (retiring its old "#,##0.0" format) and new column Q,
#    then fill in the figures.
# ---------------------------------------------------------------------------
$ws.Range("O4:O16").Copy() | Out-Null
$ws.Range("P4:P16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("Q4:Q16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("Q4").Value = 2023
$ws.Range("Q5").Value = 1781
$ws.Range("Q6").Value = 409.1
$ws.Range("Q7").Value = 435.2
$ws.Range("Q8").Value = 27.9
$ws.Range("Q9").Value = "-"
$ws.Range("Q10").Value = 194.7
$ws.Range("Q11").Value = 265.1
$ws.Range("Q12").Value = 193.9
$ws.Range("Q13").Value = 251.9
$ws.Range("Q14").Value = 1.7
$ws.Range("Q15").Value = "-"
$ws.Range("Q16").Value = 1.5

# ---------------------------------------------------------------------------
# 4) Row-height touch-ups that came along with the wider table
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 42.75
$ws.Rows(5).RowHeight = 14.25
$ws.Rows(6).RowHeight = 14.25
$ws.Rows(7).RowHeight = 14.25
$ws.Rows(8).RowHeight = 14.25
$ws.Rows(9).RowHeight = 14.25
$ws.Rows(10).RowHeight = 14.25
$ws.Rows(11).RowHeight = 14.25
$ws.Rows(16).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 5) Drop the stray cell selection left over on the sheet view
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
